$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "#firstName"
$ws.Range("C3").Value = "#lastName"

$ws.Range("C4").Value = "#userEmail"
$ws.Range("A4").Value = "Email"
$ws.Range("D4").Value = "autom@gmail.com"
$ws.Range("B4").Value = "Textbox"
$ws.Range("E4").Value = "Lname2"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:autom@gmail.com")

$ws.Columns.Item(4).ColumnWidth = 6.14

$ws.Range("A6").Select() | Out-Null
